# Generate Report for Archive
#
# The localization status moved on from "Ready for handoff" to
# "In Translation" for every tracked file, and the two status columns
# (one per sheet pair) were narrowed to fit the new, shorter label.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text shown for every tracked file.
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# Narrow the status columns to match the shorter text.
$overview.Range("E:E").ColumnWidth = 12.5
$overview.Range("F:F").ColumnWidth = 12.5
$zhcn.Range("C:C").ColumnWidth = 12.5
$dede.Range("C:C").ColumnWidth = 12.5
